$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 20 (shifts former rows 20..141 down to 21..142)
$ws.Rows.Item(20).Insert()
# Insert a second new row at 83 (shifts former row 20's-copy-now-at-83 .. down by one more)
$ws.Rows.Item(83).Insert()

# Fill in the fixed (constant across all data rows) columns for the two new rows
foreach ($r in @(20, 83)) {
    $ws.Cells.Item($r, 1).Value = 10
    $ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value = "La Araucanía"
    $ws.Cells.Item($r, 5).Value = 9
    $ws.Cells.Item($r, 6).Value = 100114007
    $ws.Cells.Item($r, 7).Value = "Jengibre"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 14).Value = "`$/caja 13 kilos"
    $ws.Cells.Item($r, 15).Value = "Perú"
    $ws.Cells.Item($r, 17).Value = 13
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

# Row 20 (new) variable values: D,I,J,K,L,M,P
$ws.Cells.Item(20, 4).Value = 44635
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 30
$ws.Cells.Item(20, 11).Value = 25000
$ws.Cells.Item(20, 12).Value = 25000
$ws.Cells.Item(20, 13).Value = 25000
$ws.Cells.Item(20, 16).Value = 1923

# Row 83 (new) variable values: D,I,J,K,L,M,P
$ws.Cells.Item(83, 4).Value = 44634
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 70
$ws.Cells.Item(83, 11).Value = 22000
$ws.Cells.Item(83, 12).Value = 25000
$ws.Cells.Item(83, 13).Value = 23714
$ws.Cells.Item(83, 16).Value = 1824
